# The "reviews_count" column (column E) is empty for every row in this
# sheet and is being removed entirely, shifting the remaining columns
# (reviews_average, latitude, longitude, is_permanently_closed,
# gmaps_link, latest_review_date) one position to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).Delete()
